# Append four new daily rows (2021-02-06 .. 2021-02-09) to the Indiana
# hospital ventilator dataset, mirroring the shape of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 348; Date = "2021-02-06"; Values = @(2175, 290, 1137, 748, 2801, 153, 442, 2206, 13.33, 52.28, 34.39, 5.46, 15.78, 78.76000000000001) },
    @{ Row = 349; Date = "2021-02-07"; Values = @(2173, 285, 1079, 809, 2797, 146, 422, 2229, 13.12, 49.65, 37.23, 5.22, 15.09, 79.69) },
    @{ Row = 350; Date = "2021-02-08"; Values = @(2172, 274, 1081, 817, 2798, 142, 430, 2226, 12.62, 49.77, 37.62, 5.08, 15.37, 79.56) },
    @{ Row = 351; Date = "2021-02-09"; Values = @(2174, 262, 1141, 771, 2797, 138, 463, 2196, 12.05, 52.48, 35.46, 4.93, 16.55, 78.51000000000001) }
)

foreach ($item in $newRows) {
    $r = $item.Row

    # Write the DATE column as a literal text string (not an auto-converted
    # date serial). Using a formula that evaluates to a text literal and then
    # collapsing it to a static value via copy/paste-values keeps the cell a
    # plain shared-string with the default (unstyled) formatting, just like
    # the rest of column A.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Formula = '="' + $item.Date + '"'
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)  # xlPasteValues

    # Columns B..O hold the numeric measurements.
    $col = 2
    foreach ($v in $item.Values) {
        $ws.Cells.Item($r, $col).Value = $v
        $col++
    }
}

$excel.CutCopyMode = 0

Write-Host ("A348=" + $ws.Range("A348").Value() + " O348=" + $ws.Range("O348").Value())
Write-Host ("A351=" + $ws.Range("A351").Value() + " O351=" + $ws.Range("O351").Value())
